# Update the acquisition timestamp column (A2:A10) on the "ランサーズ" sheet
# from "2025-09-13 18:21:56" to "2025-09-13 18:28:17".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-13 18:28:17"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
